$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title 1: "A" + " " + "slide" -> "A slide"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "A slide"

# Table cell (row 1, col 2): "a" + " " + "table" -> "a table"
$s.Shapes.Item(3).Table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "a table"

# TextBox 3: "Plus" + " " + "an" + " " + "image" -> "Plus an image"
$s.Shapes.Item(7).TextFrame.TextRange.Text = "Plus an image"
